# "Add logs REST 1000"
# Adds a new "REST 3-3-1" logs series (columns E/F) to the "100" and "1000"
# sheets, wires up their averages, updates the summary numbers on Sheet1,
# and nudges a couple of view/selection + column-width cosmetics to match.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # "Sheet1"
$ws2 = $wb.Worksheets.Item(2)   # "1"
$ws4 = $wb.Worksheets.Item(4)   # "100"
$ws5 = $wb.Worksheets.Item(5)   # "1000"

# ---------------------------------------------------------------------------
# Sheet1 : summary table updates
# ---------------------------------------------------------------------------

# Row 13 (REST, 1 req) - refine I/J values
$ws1.Range("I13").Value = 138.69764079999999
$ws1.Range("J13").Value = 176.93045740000002

# Row 15 (REST, 100 req) - C:F become "-" placeholders, G/H swap to the
# new precise values, I/J pick up the new REST-3-3-1 averages
$ws1.Range("C15").Value = "-"
$ws1.Range("D15").Value = "-"
$ws1.Range("E15").Value = "-"
$ws1.Range("F15").Value = "-"
$ws1.Range("G15").Value = 13.750628653488187
$ws1.Range("H15").Value = 15.730368793215744
$ws1.Range("I15").Value = 4.3084736182900008
$ws1.Range("J15").Value = 6.1620120192100023

# Row 16 (REST, 1000 req) - E/F become "-" placeholders too, G/H refined,
# I/J new REST-3-3-1 averages
$ws1.Range("E16").Value = "-"
$ws1.Range("F16").Value = "-"
$ws1.Range("G16").Value = 37.602944067000003
$ws1.Range("H16").Value = 43.573344314000003
$ws1.Range("I16").Value = 17.184743181666665
$ws1.Range("J16").Value = 21.656012895

# Column J (10th column) got a bit wider
$ws1.Columns.Item(10).ColumnWidth = 12.5

# Selection cosmetics
$ws1.Activate()
$ws1.Range("I26").Select()

# ---------------------------------------------------------------------------
# "1" sheet : just a selection change
# ---------------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("E13").Select()

# ---------------------------------------------------------------------------
# "100" sheet : new REST 3-3-1 log columns (E/F)
# ---------------------------------------------------------------------------
$ws4.Range("E1").Value = "eks 3-3-1"

$ws4.Range("E2").Value = 2.7292020950000002
$ws4.Range("F2").Value = 2.6596028999999999
$ws4.Range("E3").Value = 2.9139598539999998
$ws4.Range("F3").Value = 2.7710083000000001
$ws4.Range("E4").Value = 3.1439539270000001
$ws4.Range("F4").Value = 2.4935236999999999
$ws4.Range("E5").Value = 3.1168107580000002
$ws4.Range("F5").Value = 2.9569702000000002
$ws4.Range("E6").Value = 3.6095598780000002
$ws4.Range("F6").Value = 2.9653646
$ws4.Range("E7").Value = 3.5372892679999999
$ws4.Range("F7").Value = 3.1649807999999999
$ws4.Range("E8").Value = 4.5987032460000004
$ws4.Range("F8").Value = 3.1833494999999998
$ws4.Range("E9").Value = 4.4964910839999996
$ws4.Range("F9").Value = 3.2757296
$ws4.Range("E10").Value = 4.8271058739999999
$ws4.Range("F10").Value = 3.8505536
$ws4.Range("E11").Value = 5.2292365260000002
$ws4.Range("F11").Value = 3.4217824000000001
$ws4.Range("E12").Value = 5.1695034389999996
$ws4.Range("F12").Value = 3.4980422
$ws4.Range("E13").Value = 5.5189861010000003
$ws4.Range("F13").Value = 3.5828985000000002
$ws4.Range("E14").Value = 5.3165870320000002
$ws4.Range("F14").Value = 3.5399067999999998
$ws4.Range("E15").Value = 5.5389775529999996
$ws4.Range("F15").Value = 3.2353968000000002
$ws4.Range("E16").Value = 5.4808187220000004
$ws4.Range("F16").Value = 3.5354155
$ws4.Range("E17").Value = 5.4470838060000002
$ws4.Range("F17").Value = 3.1915398000000001
$ws4.Range("E18").Value = 5.585373594
$ws4.Range("F18").Value = 3.3686238999999998
$ws4.Range("E19").Value = 5.3054562809999997
$ws4.Range("F19").Value = 3.5266548000000002
$ws4.Range("E20").Value = 5.4838470340000001
$ws4.Range("F20").Value = 3.4087679
$ws4.Range("E21").Value = 5.3643920979999997
$ws4.Range("F21").Value = 3.4928642000000001
$ws4.Range("E22").Value = 5.564598814
$ws4.Range("F22").Value = 3.8123749
$ws4.Range("E23").Value = 5.576904334
$ws4.Range("F23").Value = 4.1859776999999996
$ws4.Range("E24").Value = 5.8230850920000004
$ws4.Range("F24").Value = 3.2420076
$ws4.Range("E25").Value = 5.7911737499999996
$ws4.Range("F25").Value = 3.9184050099999999
$ws4.Range("E26").Value = 6.0380679900000001
$ws4.Range("F26").Value = 4.10686252
$ws4.Range("E27").Value = 6.0833704510000004
$ws4.Range("F27").Value = 3.6357332000000002
$ws4.Range("E28").Value = 5.8982805999999997
$ws4.Range("F28").Value = 3.6813950520000001
$ws4.Range("E29").Value = 5.9941151819999998
$ws4.Range("F29").Value = 4.6583094000000003
$ws4.Range("E30").Value = 6.0378043789999998
$ws4.Range("F30").Value = 3.3637171000000001
$ws4.Range("E31").Value = 6.1409867250000003
$ws4.Range("F31").Value = 3.9316867000000002
$ws4.Range("E32").Value = 6.1179224650000004
$ws4.Range("F32").Value = 3.4166753999999999
$ws4.Range("E33").Value = 6.2360980110000002
$ws4.Range("F33").Value = 3.4902481999999999
$ws4.Range("E34").Value = 5.8665021900000003
$ws4.Range("F34").Value = 3.5571111800000002
$ws4.Range("E35").Value = 6.2284837020000001
$ws4.Range("F35").Value = 4.8283408000000003
$ws4.Range("E36").Value = 6.3322860570000001
$ws4.Range("F36").Value = 3.6437355999999999
$ws4.Range("E37").Value = 6.2266006489999999
$ws4.Range("F37").Value = 3.7052955000000001
$ws4.Range("E38").Value = 6.560785804
$ws4.Range("F38").Value = 3.9050850000000001
$ws4.Range("E39").Value = 6.2869866639999996
$ws4.Range("F39").Value = 4.6746353000000003
$ws4.Range("E40").Value = 6.212184014
$ws4.Range("F40").Value = 3.9259944
$ws4.Range("E41").Value = 6.519280416
$ws4.Range("F41").Value = 4.8478918000000002
$ws4.Range("E42").Value = 6.4699012099999997
$ws4.Range("F42").Value = 4.1798651199999997
$ws4.Range("E43").Value = 6.1531579550000002
$ws4.Range("F43").Value = 4.3754077999999996
$ws4.Range("E44").Value = 6.3862713260000001
$ws4.Range("F44").Value = 4.7767412
$ws4.Range("E45").Value = 6.2192567480000003
$ws4.Range("F45").Value = 3.8506687999999998
$ws4.Range("E46").Value = 6.2869239700000001
$ws4.Range("F46").Value = 4.5060402599999998
$ws4.Range("E47").Value = 6.5656919519999999
$ws4.Range("F47").Value = 4.4578930000000003
$ws4.Range("E48").Value = 6.2167255969999999
$ws4.Range("F48").Value = 4.0761202000000001
$ws4.Range("E49").Value = 6.4600795460000002
$ws4.Range("F49").Value = 3.3852497000000001
$ws4.Range("E50").Value = 6.473771878
$ws4.Range("F50").Value = 3.9930544000000001
$ws4.Range("E51").Value = 6.6505234980000001
$ws4.Range("F51").Value = 3.9562780000000002
$ws4.Range("E52").Value = 6.6377352509999996
$ws4.Range("F52").Value = 4.4237745000000004
$ws4.Range("E53").Value = 6.3115777800000004
$ws4.Range("F53").Value = 4.4766726099999996
$ws4.Range("E54").Value = 6.3890873160000003
$ws4.Range("F54").Value = 3.9798605
$ws4.Range("E55").Value = 6.2415698959999997
$ws4.Range("F55").Value = 3.9136272999999999
$ws4.Range("E56").Value = 6.2071788630000002
$ws4.Range("F56").Value = 3.6371042
$ws4.Range("E57").Value = 6.6794662809999998
$ws4.Range("F57").Value = 5.3780419000000004
$ws4.Range("E58").Value = 6.7095222310000002
$ws4.Range("F58").Value = 4.4100102000000003
$ws4.Range("E59").Value = 6.7302938220000001
$ws4.Range("F59").Value = 4.6377058
$ws4.Range("E60").Value = 6.6177953819999997
$ws4.Range("F60").Value = 5.3390423
$ws4.Range("E61").Value = 6.4228901389999997
$ws4.Range("F61").Value = 4.6269254999999996
$ws4.Range("E62").Value = 6.6581084649999998
$ws4.Range("F62").Value = 4.0187423000000004
$ws4.Range("E63").Value = 6.7782667529999996
$ws4.Range("F63").Value = 4.3591581000000001
$ws4.Range("E64").Value = 6.6298860089999998
$ws4.Range("F64").Value = 5.3143018
$ws4.Range("E65").Value = 6.4391638999999996
$ws4.Range("F65").Value = 4.2059369770000004
$ws4.Range("E66").Value = 6.8219661360000003
$ws4.Range("F66").Value = 4.5312963999999996
$ws4.Range("E67").Value = 6.697374194
$ws4.Range("F67").Value = 5.3757216999999997
$ws4.Range("E68").Value = 6.7323871659999996
$ws4.Range("F68").Value = 5.1588295999999998
$ws4.Range("E69").Value = 6.8254081480000002
$ws4.Range("F69").Value = 4.8623154
$ws4.Range("E70").Value = 6.8129000839999998
$ws4.Range("F70").Value = 4.9291479000000002
$ws4.Range("E71").Value = 6.6607811410000002
$ws4.Range("F71").Value = 4.3111549
$ws4.Range("E72").Value = 6.7392139579999997
$ws4.Range("F72").Value = 5.1220546000000002
$ws4.Range("E73").Value = 6.6252405019999996
$ws4.Range("F73").Value = 5.3133235000000001
$ws4.Range("E74").Value = 6.7809710550000002
$ws4.Range("F74").Value = 4.6966323000000001
$ws4.Range("E75").Value = 6.8740536219999999
$ws4.Range("F75").Value = 5.2325647999999996
$ws4.Range("E76").Value = 6.8191149170000003
$ws4.Range("F76").Value = 4.1134607000000001
$ws4.Range("E77").Value = 6.9005769050000003
$ws4.Range("F77").Value = 5.5676473
$ws4.Range("E78").Value = 6.996686167
$ws4.Range("F78").Value = 4.9003757999999999
$ws4.Range("E79").Value = 6.7046614670000002
$ws4.Range("F79").Value = 5.1176437999999997
$ws4.Range("E80").Value = 6.8529256419999998
$ws4.Range("F80").Value = 5.1366923
$ws4.Range("E81").Value = 6.9999965639999999
$ws4.Range("F81").Value = 5.1673888000000003
$ws4.Range("E82").Value = 6.9669746080000001
$ws4.Range("F82").Value = 5.4927469999999996
$ws4.Range("E83").Value = 6.655363973
$ws4.Range("F83").Value = 5.0173715999999997
$ws4.Range("E84").Value = 7.0528115800000002
$ws4.Range("F84").Value = 4.9957159000000004
$ws4.Range("E85").Value = 6.7012618760000002
$ws4.Range("F85").Value = 4.4617491999999999
$ws4.Range("E86").Value = 6.8289866049999999
$ws4.Range("F86").Value = 5.3272054000000004
$ws4.Range("E87").Value = 6.7347660549999997
$ws4.Range("F87").Value = 5.2625881999999997
$ws4.Range("E88").Value = 6.8580385939999999
$ws4.Range("F88").Value = 5.3740496000000002
$ws4.Range("E89").Value = 6.8045120710000004
$ws4.Range("F89").Value = 4.9401090999999999
$ws4.Range("E90").Value = 6.9116560070000004
$ws4.Range("F90").Value = 5.1381157000000002
$ws4.Range("E91").Value = 6.867791134
$ws4.Range("F91").Value = 4.6208998000000001
$ws4.Range("E92").Value = 6.834486751
$ws4.Range("F92").Value = 4.8576239000000001
$ws4.Range("E93").Value = 6.9510964380000004
$ws4.Range("F93").Value = 5.3338840000000003
$ws4.Range("E94").Value = 6.8751466639999999
$ws4.Range("F94").Value = 5.1581833000000001
$ws4.Range("E95").Value = 7.0154446520000002
$ws4.Range("F95").Value = 5.3624752000000004
$ws4.Range("E96").Value = 6.9260482249999997
$ws4.Range("F96").Value = 5.2062884
$ws4.Range("E97").Value = 7.1113083750000001
$ws4.Range("F97").Value = 4.7772959999999998
$ws4.Range("E98").Value = 6.8360576020000003
$ws4.Range("F98").Value = 5.0907201999999998
$ws4.Range("E99").Value = 7.0772938749999996
$ws4.Range("F99").Value = 5.8177897999999999
$ws4.Range("E100").Value = 6.8878230739999999
$ws4.Range("F100").Value = 5.3739559000000003
$ws4.Range("E101").Value = 7.1763748669999998
$ws4.Range("F101").Value = 5.7696630000000004

$ws4.Range("E103").Formula = "=AVERAGE(E2:E101)"
$ws4.Range("F103").Formula = "=AVERAGE(F2:F101)"

# New column widths for D/E/F
$ws4.Columns.Item(4).ColumnWidth = 17.666666666666668
$ws4.Columns.Item(5).ColumnWidth = 15.666666666666666
$ws4.Columns.Item(6).ColumnWidth = 16.166666666666668

# View/selection cosmetics
$ws4.Activate()
$excel.ActiveWindow.ScrollRow = 615
$excel.ActiveWindow.ScrollColumn = 1
$ws4.Range("D640").Select()

# ---------------------------------------------------------------------------
# "1000" sheet : new REST 3-3-1 log columns (E/F)
# ---------------------------------------------------------------------------
$ws5.Range("E1").Value = "REST 3-3-1"

$ws5.Range("E2").Value = 20.492216899999999
$ws5.Range("F2").Value = 16.025558499999999
$ws5.Range("E3").Value = 21.34514136
$ws5.Range("F3").Value = 16.244249799999999
$ws5.Range("E4").Value = 21.37176861
$ws5.Range("F4").Value = 15.7957144
$ws5.Range("E5").Value = 21.800629149999999
$ws5.Range("F5").Value = 20.504618700000002
$ws5.Range("E6").Value = 22.26721723
$ws5.Range("F6").Value = 17.2695373
$ws5.Range("E7").Value = 22.659104119999999
$ws5.Range("F7").Value = 17.26878039

$ws5.Range("E9").Formula = "=AVERAGE(E2:E7)"
$ws5.Range("F9").Formula = "=AVERAGE(F2:F7)"

# New column widths for E/F
$ws5.Columns.Item(5).ColumnWidth = 18.0
$ws5.Columns.Item(6).ColumnWidth = 11.833333333333332

# View/selection cosmetics
$ws5.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws5.Range("F9").Select()

# Leave Sheet1 as the active sheet, matching tabSelected="1" in sheet1.xml
$ws1.Activate()
